$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update column D (Website / Source) with real URLs instead of
#        repeating the product name. Columns A/B/C are unchanged. ---
$ws.Range("D3").Value  = "https://github.com/features/copilot"
$ws.Range("D4").Value  = "https://www.tabnine.com/"
$ws.Range("D5").Value  = "https://aws.amazon.com/q/developer/"
$ws.Range("D6").Value  = "https://sourcegraph.com/cody"
$ws.Range("D7").Value  = "https://windsurf.com/"
$ws.Range("D8").Value  = "https://mutable.ai/"
$ws.Range("D9").Value  = "https://safurai.com/"
$ws.Range("D10").Value = "https://replit.com/ai"
$ws.Range("D11").Value = "https://www.askcodi.com/"
$ws.Range("D12").Value = "https://www.qodo.ai/"
$ws.Range("D14").Value = "https://github.com/Significant-Gravitas/AutoGPT"
$ws.Range("D15").Value = "https://cursor.com/en"
$ws.Range("D17").Value = "https://taskmatrix.ai/lander"
$ws.Range("D18").Value = "https://github.com/OpenBMB/ChatDev"
$ws.Range("D19").Value = "https://github.com/OpenBMB/ChatDev"
$ws.Range("D20").Value = "https://github.com/ali-ce/gptcode-ui"
$ws.Range("D21").Value = "https://athenian.com/"
$ws.Range("D22").Value = "https://www.raycast.com/"
# D13 stays "-" (no source link) and D16 already holds its URL text.

# --- 2. Rebuild the hyperlinks so each one's target matches the new URL
#        text above (old ones pointed at stale / shortened URLs). The
#        underlying store only lets us add hyperlinks cleanly, so drop the
#        whole collection first and re-add in the original order. ---
$ws.Hyperlinks.Delete()

$targets = [ordered]@{
    "D3"  = "https://github.com/features/copilot"
    "D4"  = "https://www.tabnine.com/"
    "D5"  = "https://aws.amazon.com/q/developer/"
    "D6"  = "https://sourcegraph.com/cody"
    "D7"  = "https://windsurf.com/"
    "D8"  = "https://mutable.ai/"
    "D9"  = "https://safurai.com/"
    "D10" = "https://replit.com/ai"
    "D11" = "https://www.askcodi.com/"
    "D12" = "https://www.qodo.ai/"
    "D14" = "https://github.com/Significant-Gravitas/AutoGPT"
    "D15" = "https://cursor.com/en"
    "D17" = "https://taskmatrix.ai/lander"
    "D19" = "https://github.com/OpenBMB/ChatDev"
    "D20" = "https://github.com/ali-ce/gptcode-ui"
    "D21" = "https://athenian.com/"
    "D22" = "https://www.raycast.com/"
    "D18" = "https://github.com/OpenBMB/ChatDev"
}

foreach ($ref in $targets.Keys) {
    $cell = $ws.Range($ref)
    $ws.Hyperlinks.Add($cell, $targets[$ref])
    # Adding a hyperlink re-applies the built-in "Hyperlink" cell style,
    # which drops the wrap-text alignment these cells already had. Restore
    # it so formatting matches what was there before.
    $cell.WrapText = $true
    $cell.HorizontalAlignment = -4131
    $cell.VerticalAlignment = -4108
}

# --- 3. Misc view state that Excel recorded for this edit. ---
$ws.Range("C24").Select()
